# Covid_19_Dataset_and_References/References/13.xlsx
#
# The commit "cleans" the Authors column (E2:E33) by writing a refreshed
# copy of each author-list string. The refreshed copies are identical to
# the originals except that every run of two-or-more spaces (the padding
# between "%<flag>," and the next author's name) gains one extra space
# character; single inter-word spaces (e.g. "Jennifer M") are untouched.
#
# Rather than retype 32 long, accent-laden strings by hand (risking
# transcription errors), read each existing cell's value back out of the
# workbook and apply the same whitespace-padding transform, then write the
# result back into the same cell.

function Add-PadSpace($s) {
    $new = ""
    $i = 0
    $len = $s.Length
    while ($i -lt $len) {
        $ch = $s.Substring($i, 1)
        if ($ch -eq ' ') {
            $j = $i
            while ($j -lt $len -and $s.Substring($j, 1) -eq ' ') {
                $j = $j + 1
            }
            $runlen = $j - $i
            if ($runlen -gt 1) {
                $runlen = $runlen + 1
            }
            $rep = ""
            for ($k = 0; $k -lt $runlen; $k++) {
                $rep = $rep + " "
            }
            $new = $new + $rep
            $i = $j
        } else {
            $new = $new + $ch
            $i = $i + 1
        }
    }
    return $new
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2..33 (row 1 is the header "Authors" column E).
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $old = $cell.Value()
    $new = Add-PadSpace $old
    $cell.Value = $new
}
